$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1550
$ws.Range("J43").Value = 1733.3334
$ws.Range("L43").Value = 1733.3334
$ws.Range("N43").Value = -1871.3334
$ws.Range("H53").Value = 615.25
$ws.Range("I53").Value = 871.4
$ws.Range("K53").Value = 871.4
$ws.Range("M53").Value = -234.4
$ws.Range("H62").Value = 7400
$ws.Range("I62").Value = 2800
$ws.Range("K62").Value = 2800
$ws.Range("M62").Value = -2176
$ws.Range("H65").Value = 7400
$ws.Range("I65").Value = 2800
$ws.Range("K65").Value = 14000
$ws.Range("M65").Value = -10880
$ws.Range("H116").Value = 10710
$ws.Range("I116").Value = 10798.75
$ws.Range("K116").Value = 10798.75
$ws.Range("M116").Value = -7356.75
$ws.Range("H137").Value = 762770.4
$ws.Range("I137").Value = 1902026.4
$ws.Range("K137").Value = 5706079.199999999
$ws.Range("M137").Value = -5703529.199999999
$ws.Range("H141").Value = 2224.9
$ws.Range("I141").Value = 1843.75
$ws.Range("J141").Value = 3749.5
$ws.Range("K141").Value = 5531.25
$ws.Range("L141").Value = 11248.5
$ws.Range("M141").Value = -351.25
$ws.Range("N141").Value = -21608.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3176.6
$ws.Range("I45").Value = 2630
$ws.Range("J45").Value = 3996.5
$ws.Range("K45").Value = 2630
$ws.Range("L45").Value = 3996.5
$ws.Range("M45").Value = -2253
$ws.Range("N45").Value = -4750.5
$ws.Range("H61").Value = 3818.1875
$ws.Range("I61").Value = 2509.2
$ws.Range("J61").Value = 5999.8335
$ws.Range("K61").Value = 2509.2
$ws.Range("L61").Value = 5999.8335
$ws.Range("M61").Value = -2297.2
$ws.Range("N61").Value = -6423.8335
$ws.Range("H136").Value = 3818.1875
$ws.Range("I136").Value = 2509.2
$ws.Range("K136").Value = 7527.599999999999
$ws.Range("L136").Value = 17999.5005
$ws.Range("M136").Value = -4977.599999999999
$ws.Range("N136").Value = -23099.5005

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 11548113
$ws.Range("I7").Value = 12222274
$ws.Range("J7").Value = 10031250
$ws.Range("K7").Value = 12222274
$ws.Range("L7").Value = 10031250
$ws.Range("M7").Value = -12222161
$ws.Range("N7").Value = -10031476
$ws.Range("H86").Value = 5827.6665
$ws.Range("I86").Value = 1655.3334
$ws.Range("K86").Value = 1655.3334
$ws.Range("M86").Value = -532.3334
$ws.Range("H89").Value = 5827.6665
$ws.Range("I89").Value = 1655.3334
$ws.Range("K89").Value = 8276.666999999999
$ws.Range("M89").Value = -2660.666999999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 62.615383
$ws.Range("I7").Value = 53.555557
$ws.Range("K7").Value = 53.555557
$ws.Range("M7").Value = 59.444443
$ws.Range("I12").Value = 300
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -130
$ws.Range("N12").ClearContents()
$ws.Range("H16").Value = 212.1
$ws.Range("I16").Value = 171.83333
$ws.Range("J16").Value = 272.5
$ws.Range("K16").Value = 171.83333
$ws.Range("L16").Value = 272.5
$ws.Range("M16").Value = 115.16667
$ws.Range("N16").Value = -846.5
$ws.Range("H22").Value = 2036
$ws.Range("I22").Value = 1410
$ws.Range("J22").Value = 2975
$ws.Range("K22").Value = 1410
$ws.Range("L22").Value = 2975
$ws.Range("M22").Value = -1060
$ws.Range("N22").Value = -3675
$ws.Range("H31").Value = 6241.811
$ws.Range("J31").Value = 6863.793
$ws.Range("L31").Value = 6863.793
$ws.Range("N31").Value = -7453.793
$ws.Range("H34").Value = 6241.811
$ws.Range("J34").Value = 6863.793
$ws.Range("L34").Value = 6863.793
$ws.Range("N34").Value = -7267.793
$ws.Range("H68").Value = 79382.5
$ws.Range("J68").Value = 79382.5
$ws.Range("L68").Value = 79382.5
$ws.Range("N68").Value = -80880.5
$ws.Range("H71").Value = 79382.5
$ws.Range("J71").Value = 79382.5
$ws.Range("L71").Value = 238147.5
$ws.Range("N71").Value = -245635.5
$ws.Range("H86").Value = 6335
$ws.Range("I86").Value = 6000
$ws.Range("K86").Value = 6000
$ws.Range("M86").Value = -4877
$ws.Range("H89").Value = 6335
$ws.Range("I89").Value = 6000
$ws.Range("K89").Value = 30000
$ws.Range("M89").Value = -24384
$ws.Range("H94").Value = 5689.6665
$ws.Range("J94").Value = 5689.6665
$ws.Range("L94").Value = 5689.6665
$ws.Range("N94").Value = -6591.6665
$ws.Range("H105").Value = 1372.8889
$ws.Range("I105").Value = 1294.5
$ws.Range("K105").Value = 1294.5
$ws.Range("M105").Value = 452.5
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H113").Value = 212.1
$ws.Range("I113").Value = 171.83333
$ws.Range("J113").Value = 272.5
$ws.Range("K113").Value = 171.83333
$ws.Range("L113").Value = 272.5
$ws.Range("M113").Value = 1998.16667
$ws.Range("N113").Value = -4612.5
$ws.Range("H132").Value = 2342.9546
$ws.Range("I132").Value = 2179.4707
$ws.Range("J132").Value = 2898.8
$ws.Range("K132").Value = 6538.4121
$ws.Range("L132").Value = 8696.400000000001
$ws.Range("M132").Value = -4008.4121
$ws.Range("N132").Value = -13756.4
$ws.Range("H134").Value = 2288.4546
$ws.Range("I134").Value = 1074.5333
$ws.Range("K134").Value = 3223.5999
$ws.Range("M134").Value = -688.5999000000002

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 801.1111
$ws.Range("I122").Value = 581
$ws.Range("J122").Value = 1076.25
$ws.Range("K122").Value = 5229
$ws.Range("L122").Value = 9686.25
$ws.Range("M122").Value = -2779
$ws.Range("N122").Value = -14586.25
$ws.Range("H131").Value = 1702.9333
$ws.Range("I131").Value = 1231.909
$ws.Range("J131").Value = 2998.25
$ws.Range("K131").Value = 3695.727
$ws.Range("L131").Value = 8994.75
$ws.Range("M131").Value = 1344.273
$ws.Range("N131").Value = -19074.75

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1791.3334
$ws.Range("I122").Value = 1099.2
$ws.Range("K122").Value = 3297.6
$ws.Range("M122").Value = -847.6000000000004
$ws.Range("H132").Value = 59514.5
$ws.Range("I132").Value = 114825.445
$ws.Range("J132").Value = 4203.5557
$ws.Range("K132").Value = 344476.335
$ws.Range("L132").Value = 12610.6671
$ws.Range("M132").Value = -341946.335
$ws.Range("N132").Value = -17670.6671
$ws.Range("H135").Value = 216666.67
$ws.Range("J135").Value = 216666.67
$ws.Range("L135").Value = 216666.67
$ws.Range("N135").Value = -226806.67

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2242.9167
$ws.Range("I40").Value = 2292.2727
$ws.Range("K40").Value = 2292.2727
$ws.Range("M40").Value = -2156.2727
$ws.Range("H46").Value = 923.3333
$ws.Range("J46").Value = 923.3333
$ws.Range("L46").Value = 923.3333
$ws.Range("N46").Value = -1299.3333
$ws.Range("H132").Value = 8200.857
$ws.Range("I132").Value = 7226.5
$ws.Range("J132").Value = 9500
$ws.Range("K132").Value = 21679.5
$ws.Range("L132").Value = 28500
$ws.Range("M132").Value = -19149.5
$ws.Range("N132").Value = -33560
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 30184
$ws.Range("J104").Value = 30184
$ws.Range("L104").Value = 30184
$ws.Range("N104").Value = -37172
$ws.Range("H115").Value = 70000
$ws.Range("J115").Value = 70000
$ws.Range("L115").Value = 70000
$ws.Range("N115").Value = -73134
$ws.Range("H122").Value = 3662
$ws.Range("I122").Value = 2664.625
$ws.Range("J122").Value = 5257.8
$ws.Range("K122").Value = 7993.875
$ws.Range("L122").Value = 15773.4
$ws.Range("M122").Value = -5543.875
$ws.Range("N122").Value = -20673.4
$ws.Range("H126").Value = 3932.762
$ws.Range("J126").Value = 6788.3335
$ws.Range("L126").Value = 20365.0005
$ws.Range("N126").Value = -25305.0005
$ws.Range("H132").Value = 2034.9
$ws.Range("I132").Value = 1558.3334
$ws.Range("J132").Value = 2749.75
$ws.Range("K132").Value = 4675.0002
$ws.Range("L132").Value = 8249.25
$ws.Range("M132").Value = -2145.0002
$ws.Range("N132").Value = -13309.25
$ws.Range("H136").Value = 2355.2444
$ws.Range("I136").Value = 1764.2222
$ws.Range("K136").Value = 5292.6666
$ws.Range("M136").Value = -2742.6666
